$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A values (increment each by 1)
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# Add new rows with formulas
$ws.Range("I12").Formula = "=E2+E3+E3+E5+E4+E9"
$ws.Range("I13").Formula = "=I12/6"

# Update selection to I13
$ws.Range("I13").Select()
